# Update "想去人数" (interested count) figures for the exhibitions whose
# data is duplicated in the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 30
$wsExpo.Range("F4").Value = 47
$wsExpo.Range("F5").Value = 4979
$wsExpo.Range("F6").Value = 170
$wsExpo.Range("F7").Value = 81
$wsExpo.Range("F8").Value = 288
$wsExpo.Range("F9").Value = 44

# --- Sheet "全部类型" (all types) containing the same rows ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 30
$wsAll.Range("F8").Value = 47
$wsAll.Range("F9").Value = 4979
$wsAll.Range("F10").Value = 170
$wsAll.Range("F11").Value = 81
$wsAll.Range("F13").Value = 288
$wsAll.Range("F14").Value = 44
